$wb = $excel.ActiveWorkbook

# Overview sheet: G4 holds "Latest HO Xliff Generate Date" for a40cc3b6...
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-05 10:39:25"

# zh-cn sheet: H4 = Correspond Handoff Datetime, K4 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-05 10:39:14"
$wsZhCn.Range("K4").Value = "2016-09-05 10:40:18"

# de-de sheet: H4 = Correspond Handoff Datetime (shared string with Overview!G4), K4 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-05 10:39:25"
$wsDeDe.Range("K4").Value = "2016-09-05 10:40:38"
